$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 3.5
$ws.Range("K29").Value = 10.5
$ws.Range("M29").Value = 270.5
$ws.Range("H70").Value = 4497.5454
$ws.Range("I70").Value = 3995
$ws.Range("K70").Value = 11985
$ws.Range("M70").Value = -11715
$ws.Range("H73").Value = 4497.5454
$ws.Range("I73").Value = 3995
$ws.Range("K73").Value = 11985
$ws.Range("M73").Value = -11049
$ws.Range("H96").Value = 487.27777
$ws.Range("J96").Value = 543
$ws.Range("L96").Value = 1629
$ws.Range("N96").Value = -4375
$ws.Range("H100").Value = 3824.25
$ws.Range("I100").Value = 3099.5
$ws.Range("K100").Value = 3099.5
$ws.Range("M100").Value = -2558.5
$ws.Range("H111").Value = 1200.7273
$ws.Range("I111").Value = 1170.8
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 3512.4
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -445.3999999999996
$ws.Range("N111").Value = -10634
$ws.Range("H116").Value = 1766169.5
$ws.Range("I116").Value = 12391.615
$ws.Range("K116").Value = 12391.615
$ws.Range("M116").Value = -8949.615
$ws.Range("H134").Value = 91992
$ws.Range("J134").Value = 91992
$ws.Range("L134").Value = 91992
$ws.Range("N134").Value = -102132
$ws.Range("H136").Value = 78934
$ws.Range("J136").Value = 78934
$ws.Range("L136").Value = 78934
$ws.Range("N136").Value = -89134
$ws.Range("H137").Value = 726515.0600000001
$ws.Range("I137").Value = 1418.3572
$ws.Range("K137").Value = 4255.071599999999
$ws.Range("M137").Value = -1705.071599999999
$ws.Range("H138").Value = 1694.2565
$ws.Range("I138").Value = 1297.75
$ws.Range("J138").Value = 1830.9828
$ws.Range("K138").Value = 3893.25
$ws.Range("L138").Value = 5492.9484
$ws.Range("M138").Value = 1246.75
$ws.Range("N138").Value = -15772.9484
$ws.Range("H139").Value = 53769
$ws.Range("J139").Value = 53769
$ws.Range("L139").Value = 53769
$ws.Range("N139").Value = -64049

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2404.25
$ws.Range("I61").Value = 1790.2
$ws.Range("K61").Value = 1790.2
$ws.Range("M61").Value = -1578.2
$ws.Range("H97").Value = 1791.6666
$ws.Range("J97").Value = 3250
$ws.Range("L97").Value = 3250
$ws.Range("N97").Value = -4242
$ws.Range("H132").Value = 2156.3462
$ws.Range("J132").Value = 2621
$ws.Range("L132").Value = 7863
$ws.Range("N132").Value = -12923
$ws.Range("H136").Value = 2404.25
$ws.Range("I136").Value = 1790.2
$ws.Range("K136").Value = 5370.6
$ws.Range("M136").Value = -2820.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1534.5714
$ws.Range("I86").Value = 1289.25
$ws.Range("J86").Value = 2319.6
$ws.Range("K86").Value = 1289.25
$ws.Range("L86").Value = 2319.6
$ws.Range("M86").Value = -166.25
$ws.Range("N86").Value = -4565.6
$ws.Range("H89").Value = 1534.5714
$ws.Range("I89").Value = 1289.25
$ws.Range("J89").Value = 2319.6
$ws.Range("K89").Value = 6446.25
$ws.Range("L89").Value = 11598
$ws.Range("M89").Value = -830.25
$ws.Range("N89").Value = -22830
$ws.Range("H94").Value = 858.05554
$ws.Range("I94").Value = 592.5
$ws.Range("K94").Value = 592.5
$ws.Range("M94").Value = -141.5
$ws.Range("H99").Value = 2844003.8
$ws.Range("I99").Value = 3027.1428
$ws.Range("K99").Value = 3027.1428
$ws.Range("M99").Value = -1529.1428
$ws.Range("H110").Value = 76960
$ws.Range("J110").Value = 76960
$ws.Range("L110").Value = 76960
$ws.Range("N110").Value = -85140
$ws.Range("H112").Value = 63337.145
$ws.Range("J112").Value = 63337.145
$ws.Range("L112").Value = 63337.145
$ws.Range("N112").Value = -66291.14499999999
$ws.Range("H130").Value = 100372.25
$ws.Range("J130").Value = 100372.25
$ws.Range("L130").Value = 100372.25
$ws.Range("N130").Value = -110412.25
$ws.Range("H132").Value = 101420
$ws.Range("J132").Value = 101420
$ws.Range("L132").Value = 101420
$ws.Range("N132").Value = -111540
$ws.Range("H134").Value = 3525.1538
$ws.Range("I134").Value = 2916.2856
$ws.Range("K134").Value = 8748.856800000001
$ws.Range("M134").Value = -6213.856800000001
$ws.Range("H135").Value = 73360
$ws.Range("J135").Value = 73360
$ws.Range("L135").Value = 73360
$ws.Range("N135").Value = -83500
$ws.Range("H138").Value = 99941
$ws.Range("J138").Value = 99941
$ws.Range("L138").Value = 99941
$ws.Range("N138").Value = -110221
$ws.Range("H140").Value = 124336.375
$ws.Range("J140").Value = 70670.14
$ws.Range("L140").Value = 70670.14
$ws.Range("N140").Value = -81030.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2147.5
$ws.Range("I58").Value = 1658.75
$ws.Range("K58").Value = 1658.75
$ws.Range("M58").Value = -1455.75
$ws.Range("H136").Value = 2147.5
$ws.Range("I136").Value = 1658.75
$ws.Range("K136").Value = 4976.25
$ws.Range("M136").Value = -2426.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2415.7083
$ws.Range("I132").Value = 1990
$ws.Range("J132").Value = 2454.4092
$ws.Range("K132").Value = 17910
$ws.Range("L132").Value = 22089.6828
$ws.Range("M132").Value = -15380
$ws.Range("N132").Value = -27149.6828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1997.8334
$ws.Range("I97").Value = 1696
$ws.Range("J97").Value = 2299.6667
$ws.Range("K97").Value = 1696
$ws.Range("L97").Value = 2299.6667
$ws.Range("M97").Value = -1200
$ws.Range("N97").Value = -3291.6667
$ws.Range("H122").Value = 12691.218
$ws.Range("I122").Value = 13836.789
$ws.Range("K122").Value = 41510.367
$ws.Range("M122").Value = -39060.367
$ws.Range("H132").Value = 6578.5713
$ws.Range("J132").Value = 4850
$ws.Range("L132").Value = 14550
$ws.Range("N132").Value = -19610
$ws.Range("H135").Value = 79990
$ws.Range("J135").Value = 79990
$ws.Range("L135").Value = 79990
$ws.Range("N135").Value = -90130
$ws.Range("H140").Value = 44498.75
$ws.Range("J140").Value = 41798
$ws.Range("L140").Value = 41798
$ws.Range("N140").Value = -52158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2014.6666
$ws.Range("I46").Value = 1418.6666
$ws.Range("K46").Value = 1418.6666
$ws.Range("M46").Value = -1230.6666
$ws.Range("H68").Value = 2815.8
$ws.Range("I68").Value = 2794.5
$ws.Range("J68").Value = 2830
$ws.Range("K68").Value = 2794.5
$ws.Range("L68").Value = 2830
$ws.Range("M68").Value = -2045.5
$ws.Range("N68").Value = -4328
$ws.Range("H71").Value = 2815.8
$ws.Range("I71").Value = 2794.5
$ws.Range("J71").Value = 2830
$ws.Range("K71").Value = 13972.5
$ws.Range("L71").Value = 14150
$ws.Range("M71").Value = -10228.5
$ws.Range("N71").Value = -21638
$ws.Range("H100").Value = 2388.6
$ws.Range("I100").Value = 2340.8572
$ws.Range("K100").Value = 2340.8572
$ws.Range("M100").Value = -1799.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1987.3125
$ws.Range("I96").Value = 1816.3334
$ws.Range("J96").Value = 2089.9
$ws.Range("K96").Value = 1816.3334
$ws.Range("L96").Value = 2089.9
$ws.Range("M96").Value = -443.3334
$ws.Range("N96").Value = -4835.9
$ws.Range("H122").Value = 3678.9375
$ws.Range("I122").Value = 2489.0715
$ws.Range("K122").Value = 7467.2145
$ws.Range("M122").Value = -5017.2145
$ws.Range("H132").Value = 806362.4399999999
$ws.Range("J132").Value = 5437273.5
$ws.Range("L132").Value = 16311820.5
$ws.Range("N132").Value = -16316880.5
